# Stardew Valley Win32API spec workbook update
# 1. Collision Enter/Stay/Exit functions implemented, 2. FarmScene Object edits
# -> fill the previously-blank Progress ("D"/"E") columns with "-" placeholders
#    on both the Common and Scene sheets, adjust the Scene print area, and
#    move the active selection on each sheet.

$wb = $excel.ActiveWorkbook

$xlCenter = -4108

# ---------------------------------------------------------------------------
# Common sheet: blank D/E "progress" cells get a "-" placeholder, centered
# ---------------------------------------------------------------------------
$wsCommon = $wb.Worksheets.Item("Common")

$commonCells = @(
    "D5","D6",
    "E10","E11",
    "D12","D13",
    "D14","E14",
    "E15","E16",
    "D17","E17",
    "E18","E19","E20","E21",
    "D22","E22",
    "D23","E23",
    "D24","E24",
    "E25","E26","E27","E28","E29"
)

foreach ($addr in $commonCells) {
    $cell = $wsCommon.Range($addr)
    $cell.Value = "-"
    $cell.HorizontalAlignment = $xlCenter
    $cell.WrapText = $false
}

# ---------------------------------------------------------------------------
# Scene sheet: same treatment for its blank "progress" cells
# ---------------------------------------------------------------------------
$wsScene = $wb.Worksheets.Item("Scene")

$sceneCells = @(
    "E5",
    "E14","E15","E16","E17",
    "E19",
    "E21",
    "E23","E24",
    "E26",
    "E28",
    "E30",
    "E33","E34"
)

foreach ($addr in $sceneCells) {
    $cell = $wsScene.Range($addr)
    $cell.Value = "-"
    $cell.HorizontalAlignment = $xlCenter
    $cell.WrapText = $false
}

# ---------------------------------------------------------------------------
# Scene print area shrinks by one row (36 -> 35)
# ---------------------------------------------------------------------------
$wsScene.PageSetup.PrintArea = '$A$1:$G$35'

# ---------------------------------------------------------------------------
# Move the on-screen selection / view to match where editing happened
# ---------------------------------------------------------------------------
$wsCommon.Activate()
$wsCommon.Range("L22").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ZoomSheetLayoutView = 100

$wsScene.Activate()
$wsScene.Range("J16").Select()

$wsCommon.Activate()
